$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old rows 5-19 (rows beyond the new 4-row table)
$ws.Range("A5:B19").ClearContents()

# Update column A (run) values for rows 2-4
$ws.Range("A2").Value = "sumIntensity_1"
$ws.Range("A3").Value = "sumIntensity_2"
$ws.Range("A4").Value = "sumIntensity_3"

# Update column B (treatment) values for rows 1-4
$ws.Range("B1").Value = "treatment"
$ws.Range("B2").Value = "fold1"
$ws.Range("B3").Value = "fold4"
$ws.Range("B4").Value = "fold10"

# Update selection to match the final state
$ws.Range("A2:A4").Select()
